# Add a Fade entrance animation (triggered "With Previous") for the
# "Rectangle 2" shape (spid=3, the "github link" textbox) on slide 1,
# joining the existing click-triggered animation group alongside the
# title and the two body-text paragraph animations.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Rectangle 2" shape (the hyperlinked "github link" textbox).
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Rectangle 2") {
        $shape = $candidate
    }
}

$mainSeq = $s.TimeLine.MainSequence

# msoAnimEffectFade = 10, msoAnimTriggerWithPrevious = 2
# Using "WithPrevious" attaches the new effect to the existing click
# paragraph (alongside the other three effects already on this slide)
# rather than creating a brand-new click trigger.
$effect = $mainSeq.AddEffect($shape, 10, 0, 2)
